# Update "想去人数" (F column) values for both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 64
    6  = 2866
    8  = 1866
    10 = 93
    11 = 677
    14 = 197
    15 = 88
    17 = 28
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
